# Remove Wind onshore/offshore entries from the urbs intertemporal input workbook
$wb = $excel.ActiveWorkbook

# --- Commodity sheet: delete the WindOff / WindOn rows (rows 2:3) ---
$wsCommodity = $wb.Worksheets.Item("Commodity")
$wsCommodity.Activate()
$wsCommodity.Rows("2:3").Select()
$wsCommodity.Rows("2:3").Delete()

# --- Process sheet: delete the Wind (onshore) / Wind (offshore) rows (rows 2:3) ---
$wsProcess = $wb.Worksheets.Item("Process")
$wsProcess.Activate()
$wsProcess.Rows("2:3").Select()
$wsProcess.Rows("2:3").Delete()

# --- Process-Commodity sheet: delete the Wind (onshore) / Wind (offshore) In/Out rows (rows 2:5) ---
$wsProcessCommodity = $wb.Worksheets.Item("Process-Commodity")
$wsProcessCommodity.Activate()
$wsProcessCommodity.Rows("2:5").Select()
$wsProcessCommodity.Rows("2:5").Delete()

# --- SupIm sheet: delete the EU27.WindOff / EU27.WindOn columns (columns B:C) ---
$wsSupIm = $wb.Worksheets.Item("SupIm")
$wsSupIm.Activate()
$wsSupIm.Columns("B:C").Select()
$wsSupIm.Columns("B:C").Delete()
$wsSupIm.Range("C1").Select()
